$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("SValimaki                                       ")

# Insert the new worksheet directly after the last existing sheet so the
# tab order ends up: Sheet1, "SValimaki <spaces>", "SValimaki"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "SValimaki"

# Populate row 2 (A:D) with the tracking entry
$newSheet.Range("A2").Value = "Inside"
$newSheet.Range("B2").Value = "SValimaki                                       "
$newSheet.Range("C2").Value = 43423
$newSheet.Range("D2").Value = 0.5084813385416667

# Reuse the existing date/time number formats (styles) from the source sheet
# instead of creating brand-new style entries
$sourceSheet.Range("C2").Copy()
$newSheet.Range("C2").PasteSpecial(-4122)
$sourceSheet.Range("D2").Copy()
$newSheet.Range("D2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
